$d = $word.ActiveDocument

# Heading replacements in main body
$d.Content.Find.Execute("🎯 Core Research Context", $true, $false, $false, $false, $false, $true, 1, $false, "**Goal:** Core Research Context", 2)
$d.Content.Find.Execute("📋 Research Questions Framework", $true, $false, $false, $false, $false, $true, 1, $false, "**Requirements:** Research Questions Framework", 2)
$d.Content.Find.Execute("🔍 Investigation Methodology", $true, $false, $false, $false, $false, $true, 1, $false, "**Analysis:** Investigation Methodology", 2)
$d.Content.Find.Execute("**Token Structure Optimization**: Balance between simplicity and flexibility", $true, $false, $false, $false, $false, $true, 1, $false, "**Token Structure improvement**: Balance between simplicity and flexibility", 2)
$d.Content.Find.Execute("🎯 Expected Outcomes", $true, $false, $false, $false, $false, $true, 1, $false, "**Goal:** Expected Outcomes", 2)
$d.Content.Find.Execute("**Token Architecture Proposal**: Optimized structure based on competitive analysis", $true, $false, $false, $false, $false, $true, 1, $false, "**Token Architecture Proposal**: improved structure based on competitive analysis", 2)

# Footer date update
foreach ($section in $d.Sections) {
    $footer = $section.Footers.Item(1)
    $footer.Range.Find.Execute("global-sizing-concept-research-questions | Last edited: 2025-09-07 12:53 | Page [X] of [Y]", $true, $false, $false, $false, $false, $true, 1, $false, "global-sizing-concept-research-questions | Last edited: 2025-09-12 17:37 | Page [X] of [Y]", 2)
}

$d.Save()
